$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")
$ws.Range("A2").Value = "Test"
